$d = $word.ActiveDocument

$replacements = @(
    @("91×92=8372", "68×64=4352"),
    @("84×16=1344", "47×13=611"),
    @("13×79=1027", "41×33=1353"),
    @("34×90=3060", "53×56=2968"),
    @("59×31=1829", "43×99=4257"),
    @("32×68=2176", "76×99=7524"),
    @("23×37=851", "58×82=4756"),
    @("42×85=3570", "47×56=2632"),
    @("46×40=1840", "93×51=4743"),
    @("81×41=3321", "31×61=1891"),
    @("43×73=3139", "37×57=2109"),
    @("86×67=5762", "80×48=3840"),
    @("90×11=990", "32×14=448"),
    @("36×35=1260", "38×89=3382"),
    @("80×65=5200", "74×32=2368"),
    @("97×44=4268", "95×39=3705"),
    @("37×56=2072", "36×78=2808"),
    @("42×86=3612", "66×44=2904"),
    @("88×66=5808", "36×59=2124"),
    @("45×64=2880", "17×76=1292"),
    @("70×58=4060", "83×20=1660"),
    @("84×12=1008", "48×60=2880"),
    @("65×95=6175", "54×37=1998"),
    @("97×99=9603", "24×53=1272"),
    @("86×96=8256", "54×73=3942")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
